$d = $word.ActiveDocument

# Locate the paragraph containing only "Reports" - the last fully-typed
# bullet of the existing numbered list (numId=1). The paragraph right
# after it is an empty list item (ilvl=0) that is where the new content
# needs to be typed in.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Reports`r") {
        $anchor = $p.Next()
        break
    }
}

$newXml = @'
<w:p w:rsidR="00690D6B" w:rsidRPr="00690D6B" w:rsidRDefault="00690D6B" w:rsidP="00690D6B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Executing Tests in Parallel</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Configuring Extent Reports in selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:lastRenderedPageBreak/><w:t>Taking Screen shots in selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Executing javascript in selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>POM</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Created a project and sample code</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Added compiler plugin to maven</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Added surefire plugin to maven</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Created profiles for different testng xml file</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Created Batch file for test execution</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>GRID</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Downloaded latest GRID jar file</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Created hub</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Copied json sample file and modified as required</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Regestred Node to hub</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Updated driver creation logic</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>Executed Code on Distributed Environment</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr></w:p>
'@

$anchor.Range.InsertXML($newXml)
